$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append a new review row (row 8) ---
$ws.Range("A8").Value = "com.hamxa.shaynachim"
$ws.Range("B8").Value = "bitcoin"
$ws.Range("C8").Value = "georggini2@gmail.com"
$ws.Range("D8").Value = "jorjkluni03@gmail.com"
$ws.Range("E8").Value = "27/5/2019 15:59"
$ws.Range("F8").Value = "day and night! Iplay this game all the time"
$ws.Range("G8").Value = "no"

# Row 8 should pick up the same per-column formatting as the rest of the
# table (row 7 is a representative, fully-styled data row).
$ws.Range("A7:G7").Copy()
$ws.Range("A8:G8").PasteSpecial(-4122)

# PasteSpecial only copied formats; put the real values back.
$ws.Range("A8").Value = "com.hamxa.shaynachim"
$ws.Range("B8").Value = "bitcoin"
$ws.Range("C8").Value = "georggini2@gmail.com"
$ws.Range("D8").Value = "jorjkluni03@gmail.com"
$ws.Range("E8").Value = "27/5/2019 15:59"
$ws.Range("F8").Value = "day and night! Iplay this game all the time"
$ws.Range("G8").Value = "no"

# D8 mirrors C7/D7: a mailto: hyperlink whose display text is the email.
$ws.Hyperlinks.Add($ws.Range("D8"), "mailto:jorjkluni03@gmail.com", "", "", "jorjkluni03@gmail.com")

# Adding the hyperlink re-styles the cell (blue/underline font); restore
# the plain data style shared by the rest of the email columns.
$ws.Range("C8").Copy()
$ws.Range("D8").PasteSpecial(-4122)
$ws.Range("D8").Value = "jorjkluni03@gmail.com"

# --- Selection moves from G7 to F8 ---
$ws.Range("F8").Select()
